$d = $word.ActiveDocument

# "verifikasi KA" table row - nomenclature change:
#   "Justifikasi/bukti persetujuan awal rencana usaha dan/atau kegiatan dengan PIPPIB"
#     -> "Justifikasi/bukti kesesuaian lokasi rencana usaha dan/atau kegiatan dengan PIPPIB"
#
# The document has a couple of near-duplicate "... persetujuan awal rencana
# usaha dan/atau kegiatan ..." sentences in other, unrelated table rows, so
# first pin down the *whole* PIPPIB sentence (unique in the document) to get
# an exact Range, then do the actual text swap scoped to that Range only -
# this guarantees the edit can't accidentally land on a different row.
$anchorText = "Justifikasi/bukti persetujuan awal rencana usaha dan/atau kegiatan dengan PIPPIB"

$hit = $d.Content
$hit.Find.ClearFormatting()
$hit.Find.Execute($anchorText, $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)

if ($hit.Find.Found) {
    # Narrow further to just the fragment being renamed, still scoped to
    # this one sentence's Range so the search cannot leak elsewhere.
    $scope = $d.Range($hit.Start, $hit.End)
    $scope.Find.ClearFormatting()
    $scope.Find.Execute(" persetujuan awal", $true, $false, $false, $false, `
                         $false, $true, 1, $false, "", 0)

    if ($scope.Find.Found) {
        # Toggling a character-formatting property across the text swap
        # forces Word to keep this fragment as its own run instead of
        # silently re-merging it with the (identically formatted)
        # neighbouring runs once the temporary formatting is removed again.
        $scope.Bold = 1
        $scope.Text = " kesesuaian lokasi"
        $scope.Bold = 0
    }
}
